$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking strings
# (e.g. "1.00", "7.00", "56.786.71") are preserved exactly as authored,
# matching the original inline-string cell contents.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '56.786.71'
$ws.Range('E2').Value = '  +2.70%  '
$ws.Range('D3').Value = '3.002.94'
$ws.Range('E3').Value = '  +1.63%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '510.65'
$ws.Range('E5').Value = '  +4.81%  '
$ws.Range('D6').Value = '139.09'
$ws.Range('E6').Value = '  +6.01%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E8').Value = '  +4.43%  '
$ws.Range('D9').Value = '7.48'
$ws.Range('E9').Value = '  +6.16%  '
$ws.Range('E10').Value = '  +7.79%  '
$ws.Range('D11').Value = '0.357'
$ws.Range('E11').Value = '  +2.90%  '
$ws.Range('D12').Value = '0.129'
$ws.Range('E12').Value = '  +2.23%  '
$ws.Range('D13').Value = '3.516.75'
$ws.Range('E13').Value = '  +1.68%  '
$ws.Range('D14').Value = '25.88'
$ws.Range('E14').Value = '  +5.23%  '
$ws.Range('E15').Value = '  +12.70%  '
$ws.Range('D16').Value = '56.794.21'
$ws.Range('E16').Value = '  +2.93%  '
$ws.Range('D17').Value = '3.000.89'
$ws.Range('E17').Value = '  +1.76%  '
$ws.Range('D18').Value = '5.93'
$ws.Range('E18').Value = '  +6.93%  '
$ws.Range('D19').Value = '12.51'
$ws.Range('E19').Value = '  +4.20%  '
$ws.Range('D20').Value = '7.84'
$ws.Range('E20').Value = '  +5.55%  '
$ws.Range('D21').Value = '327.19'
$ws.Range('E21').Value = '  +3.56%  '
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').Value = '0.486'
$ws.Range('E23').Value = '  +5.44%  '
$ws.Range('D24').Value = '63.31'
$ws.Range('E24').Value = '  +5.44%  '
$ws.Range('E25').Value = '  +5.75%  '
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('D27').Value = '0.0₃0907'
$ws.Range('E27').Value = '  +8.10%  '
$ws.Range('D28').Value = '6.66'
$ws.Range('E28').Value = '  +3.14%  '
$ws.Range('D29').Value = '7.00'
$ws.Range('E29').Value = '  +7.57%  '
$ws.Range('D30').Value = '1.23'
$ws.Range('E30').Value = '  +6.51%  '
$ws.Range('E31').Value = '  +7.21%  '
$ws.Range('D32').Value = '20.67'
$ws.Range('E32').Value = '  +7.16%  '
$ws.Range('D33').Value = '154.49'
$ws.Range('E33').Value = '  +3.74%  '
$ws.Range('D34').Value = '4.56'
$ws.Range('E34').Value = '  +4.55%  '
$ws.Range('D35').Value = '5.70'
$ws.Range('E35').Value = '  +1.09%  '
$ws.Range('E36').Value = '  -0.71%  '
$ws.Range('D37').Value = '0.0680'
$ws.Range('E37').Value = '  +5.61%  '
$ws.Range('D38').Value = '23.79'
$ws.Range('E38').Value = '  +2.55%  '
$ws.Range('D39').Value = '3.034.73'
$ws.Range('E39').Value = '  +1.80%  '
$ws.Range('D40').Value = '37.07'
$ws.Range('E40').Value = '  +3.05%  '
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').Value = '2.285.75'
$ws.Range('E42').Value = '  +7.71%  '
$ws.Range('D43').Value = '0.649'
$ws.Range('E43').Value = '  +3.24%  '
$ws.Range('D44').Value = '3.69'
$ws.Range('E44').Value = '  +4.70%  '
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.99%  '
$ws.Range('E46').Value = '  +3.69%  '
$ws.Range('D47').Value = '1.96'
$ws.Range('E47').Value = '  +12.77%  '
$ws.Range('D48').Value = '5.88'
$ws.Range('E48').Value = '  +6.06%  '
$ws.Range('E49').Value = '  +2.81%  '
$ws.Range('D50').Value = '19.29'
$ws.Range('E50').Value = '  +1.63%  '
$ws.Range('E51').Value = '  +5.62%  '

# Restore the original (default) cell style now that the text values are set,
# so no stray number-format style is left applied to the cells.
$ws.Range('D2:D51').Style = 'Normal'
